# Updated cryptos list data (mirrors upstream coinranking.com scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "23.239.81"
$ws.Range("E2").Value = "  +0.85%  "
# Row 3
$ws.Range("D3").Value = "1.604.87"
$ws.Range("E3").Value = "  +0.34%  "
# Row 4
$ws.Range("E4").Value = "  +0.02%  "
# Row 5
$ws.Range("E5").Value = "  +0.05%  "
# Row 6
Set-TextValue "D6" "304.62"
$ws.Range("E6").Value = "  +0.80%  "
# Row 7
$ws.Range("E7").Value = "  -0.58%  "
# Row 8
Set-TextValue "D8" "52.39"
$ws.Range("E8").Value = "  +4.93%  "
# Row 9
Set-TextValue "D9" "0.3627"
$ws.Range("E9").Value = "  -0.66%  "
# Row 10
Set-TextValue "D10" "1.274"
$ws.Range("E10").Value = "  +1.18%  "
# Row 11
Set-TextValue "D11" "0.08152"
$ws.Range("E11").Value = "  -0.04%  "
# Row 12
Set-TextValue "D12" "1.002"
$ws.Range("E12").Value = "  +0.04%  "
# Row 13
Set-TextValue "D13" "22.90"
$ws.Range("E13").Value = "  +1.67%  "
# Row 14
$ws.Range("E14").Value = "  +0.05%  "
# Row 15
Set-TextValue "D15" "7.383"
$ws.Range("E15").Value = "  +0.28%  "
# Row 16
Set-TextValue "D16" "0.00001249"
$ws.Range("E16").Value = "  +0.10%  "
# Row 17
$ws.Range("D17").Value = "1.605.34"
$ws.Range("E17").Value = "  +0.20%  "
# Row 18
Set-TextValue "D18" "93.95"
$ws.Range("E18").Value = "  +2.26%  "
# Row 19
Set-TextValue "D19" "0.06916"
$ws.Range("E19").Value = "  +1.34%  "
# Row 20
Set-TextValue "D20" "18.16"
$ws.Range("E20").Value = "  -0.41%  "
# Row 21
Set-TextValue "D21" "6.541"
$ws.Range("E21").Value = "  +0.06%  "
# Row 22
$ws.Range("E22").Value = "  +0.21%  "
# Row 23
Set-TextValue "D23" "12.92"
$ws.Range("E23").Value = "  -1.19%  "
# Row 24
$ws.Range("D24").Value = "23.221.14"
$ws.Range("E24").Value = "  +0.79%  "
# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D25" "2.451"
$ws.Range("E25").Value = "  +3.64%  "
# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D26" "3.080"
$ws.Range("E26").Value = "  +9.44%  "
# Row 27
$ws.Range("E27").Value = "  +0.53%  "
# Row 28
$ws.Range("E28").Value = "  -0.28%  "
# Row 29
Set-TextValue "D29" "5.277"
$ws.Range("E29").Value = "  +0.91%  "
# Row 30
Set-TextValue "D30" "135.19"
$ws.Range("E30").Value = "  +0.51%  "
# Row 31
Set-TextValue "D31" "2.388"
$ws.Range("E31").Value = "  +1.75%  "
# Row 32
Set-TextValue "D32" "6.746"
$ws.Range("E32").Value = "  -1.44%  "
# Row 33
$ws.Range("D33").Value = "1.780.91"
$ws.Range("E33").Value = "  +0.22%  "
# Row 34
Set-TextValue "D34" "0.9637"
$ws.Range("E34").Value = "  +0.05%  "
# Row 35
Set-TextValue "D35" "0.07488"
$ws.Range("E35").Value = "  -1.51%  "
# Row 36
Set-TextValue "D36" "10.42"
$ws.Range("E36").Value = "  +0.88%  "
# Row 37
Set-TextValue "D37" "0.02758"
$ws.Range("E37").Value = "  +1.74%  "
# Row 38
Set-TextValue "D38" "0.2516"
$ws.Range("E38").Value = "  -0.49%  "
# Row 39
Set-TextValue "D39" "6.120"
$ws.Range("E39").Value = "  -2.16%  "
# Row 40
Set-TextValue "D40" "0.08793"
$ws.Range("E40").Value = "  -0.69%  "
# Row 41
Set-TextValue "D41" "1.416"
$ws.Range("E41").Value = "  +3.45%  "
# Row 42
Set-TextValue "D42" "0.7092"
$ws.Range("E42").Value = "  +0.59%  "
# Row 43
Set-TextValue "D43" "12.47"
$ws.Range("E43").Value = "  +0.37%  "
# Row 44
Set-TextValue "D44" "15.84"
$ws.Range("E44").Value = "  +3.91%  "
# Row 45
$ws.Range("E45").Value = "  -1.49%  "
# Row 46
Set-TextValue "D46" "2.330"
$ws.Range("E46").Value = "  +1.72%  "
# Row 47
$ws.Range("E47").Value = "  +0.36%  "
# Row 48
Set-TextValue "D48" "133.97"
$ws.Range("E48").Value = "  +1.42%  "
# Row 49
Set-TextValue "D49" "0.07943"
$ws.Range("E49").Value = "  +0.49%  "
# Row 50
Set-TextValue "D50" "1.208"
$ws.Range("E50").Value = "  -0.84%  "
# Row 51
Set-TextValue "D51" "1.192"
$ws.Range("E51").Value = "  -3.19%  "

Write-Output "Updated cryptos list"
